# Weekly update: insert two new daily price rows (date 2022-03-08 / serial 44628)
# for "Poroto granado" at Vega Central Mapocho de Santiago, pushing the prior
# two most-recent rows (253, 254) down to (255, 256).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current last two data rows (253, 254),
# shifting the existing rows 253/254 down to 255/256.
$ws.Rows.Item(253).Insert()
$ws.Rows.Item(253).Insert()

# New row 253: Region Metropolitana entry for 2022-03-08
$ws.Cells.Item(253, 1).Value = 9
$ws.Cells.Item(253, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(253, 3).Value = "Metropolitana"
$ws.Cells.Item(253, 4).Value = 44628
$ws.Cells.Item(253, 5).Value = 13
$ws.Cells.Item(253, 6).Value = 100112030
$ws.Cells.Item(253, 7).Value = "Poroto granado"
$ws.Cells.Item(253, 8).Value = "Sin especificar"
$ws.Cells.Item(253, 9).Value = "Primera"
$ws.Cells.Item(253, 10).Value = 52
$ws.Cells.Item(253, 11).Value = 20000
$ws.Cells.Item(253, 12).Value = 22000
$ws.Cells.Item(253, 13).Value = 21000
$ws.Cells.Item(253, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(253, 15).Value = "Región Metropolitana"
$ws.Cells.Item(253, 16).Value = 840
$ws.Cells.Item(253, 17).Value = 25
$ws.Cells.Item(253, 18).Value = "Hortaliza"

# New row 254: Region de O'Higgins entry for 2022-03-08
$ws.Cells.Item(254, 1).Value = 9
$ws.Cells.Item(254, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(254, 3).Value = "Metropolitana"
$ws.Cells.Item(254, 4).Value = 44628
$ws.Cells.Item(254, 5).Value = 13
$ws.Cells.Item(254, 6).Value = 100112030
$ws.Cells.Item(254, 7).Value = "Poroto granado"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 61
$ws.Cells.Item(254, 11).Value = 22000
$ws.Cells.Item(254, 12).Value = 23000
$ws.Cells.Item(254, 13).Value = 22492
$ws.Cells.Item(254, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(254, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(254, 16).Value = 900
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
